$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("transitions")

# ---------------------------------------------------------------------------
# Small helper: paste the formatting currently sitting on the clipboard onto
# a list of single-cell addresses (a multi-area Range only "takes" on its
# first area for both .Value assignment and PasteSpecial, so cells are
# targeted one at a time).
# ---------------------------------------------------------------------------
function Set-CellsFormat($cells) {
    foreach ($addr in $cells) {
        $ws.Range($addr).PasteSpecial(-4122) | Out-Null
    }
}

function Set-CellsValue($cells, $value) {
    foreach ($addr in $cells) {
        $ws.Range($addr).Value = $value
    }
}

# ---------------------------------------------------------------------------
# 1) Flip a handful of existing FSM transitions from T_ERROR (or T_URL) over
#    to the (new) S_IPV6 / S_URL states. These happen at their original row
#    numbers, which all sit above the row we're about to insert, so the
#    insert below doesn't disturb them.
# ---------------------------------------------------------------------------
$existingEdits = @(
    "U25",
    "S28", "T28", "U28", "V28",
    "S29", "U29",
    "O33", "S33", "T33", "U33", "V33",
    "S34", "T34", "U34", "V34",
    "O52",
    "F56", "G56"
)

# All of the above land on the "plain state name" look (regular, non-bold) --
# match the formatting already used for that look elsewhere (reference: D25,
# the S_NUMBER state-name cell).
$ws.Range("D25").Copy() | Out-Null
Set-CellsFormat $existingEdits

Set-CellsValue @("U25") "S_IPV6"
Set-CellsValue @("S28", "T28", "U28", "V28") "S_IPV6"
Set-CellsValue @("S29", "U29") "S_IPV6"
Set-CellsValue @("O33", "S33", "T33", "U33", "V33") "S_IPV6"
Set-CellsValue @("S34", "T34", "U34", "V34") "S_IPV6"
Set-CellsValue @("O52") "S_IPV6"
Set-CellsValue @("F56", "G56") "S_URL"

# ---------------------------------------------------------------------------
# 2) Insert a brand-new FSM row for the S_IPV6 state right above the old
#    row 59 (S_EQUAL), shifting everything from row 59 down by one.
# ---------------------------------------------------------------------------
$ws.Rows.Item(59).Insert()

# Columns that loop back into S_IPV6 (bare class chars: digits, colon, hex
# letters/"h"/"e"/dot) -- regular look (reference: D25, the S_NUMBER cell).
$sIpv6Cols = @("D59", "E59", "O59", "S59", "T59", "U59", "V59", "AF59")

# Columns that finalize into the T_IPV6 token -- bold look (reference: B25,
# the T_INTEGER cell).
$tIpv6Cols = @("B59", "C59", "F59", "G59", "H59", "I59", "J59", "K59", "L59", "Y59", "AD59", "AK59", "AO59")

# Remaining columns (besides A59) are illegal-transition/T_ERROR -- bold red
# look (reference: X25, the T_ERROR cell).
$tErrorCols = @("M59", "N59", "P59", "Q59", "R59", "W59", "X59", "Z59", "AA59", "AB59", "AC59", "AE59", "AG59", "AH59", "AI59", "AJ59", "AL59", "AM59", "AN59")

# Column A (state-name column) gets the bold + right-border look used by
# every other state-name cell in column A (reference: A25).
$ws.Range("A25").Copy() | Out-Null
Set-CellsFormat @("A59")

$ws.Range("D25").Copy() | Out-Null
Set-CellsFormat $sIpv6Cols

$ws.Range("B25").Copy() | Out-Null
Set-CellsFormat $tIpv6Cols

$ws.Range("X25").Copy() | Out-Null
Set-CellsFormat $tErrorCols

# Now drop in the actual text for every cell in the new row.
Set-CellsValue @("A59") "S_IPV6"
Set-CellsValue $sIpv6Cols "S_IPV6"
Set-CellsValue $tIpv6Cols "T_IPV6"
Set-CellsValue $tErrorCols "T_ERROR"

# ---------------------------------------------------------------------------
# 3) View/selection bookkeeping: "transitions" becomes the active sheet/tab,
#    with the cursor parked on H32; "float transitions" no longer is.
# ---------------------------------------------------------------------------
$ws.Activate()
$ws.Range("H32").Select()
